# Insert a new row "Clear console" / "Control + L" right above the
# existing "Clear workspace" row (currently row 7), pushing all the
# rows below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 7 (shifts "Clear workspace" etc. down to row 8+).
$ws.Rows("7:7").Insert()

# The inserted row inherits formatting from the row above (the bold
# "Description"/"Syntax" header row) - clear that so the new row looks
# like the other plain data rows.
$ws.Range("A7:B7").ClearFormats()

# Populate the new row.
$ws.Range("A7").Value = "Clear console"
$ws.Range("B7").Value = "Control + L"

# Update the active selection to match where the author last left off.
$ws.Range("B8").Select() | Out-Null
